$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"0.040495"
$ws.Range("H2").Value = [double]"0.121485"
$ws.Range("I2").Value = [double]"0.002191743187342868"
$ws.Range("J2").Value = [double]"0.002191743187342869"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"0.040495"
$ws.Range("N2").Value = [double]"0.121485"
$ws.Range("O2").Value = [double]"0.002191743187342868"
$ws.Range("P2").Value = [double]"0.002191743187342869"
$ws.Range("Q2").Value = [double]"0.001639845025"
$ws.Range("R2").Value = [double]"0.014758605225"
$ws.Range("S2").Value = [double]"4.803738199263875e-06"
$ws.Range("T2").Value = [double]"4.803738199263877e-06"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"0.040495"
$ws.Range("H3").Value = [double]"0.121485"
$ws.Range("I3").Value = [double]"0.002191743187342868"
$ws.Range("J3").Value = [double]"0.002191743187342869"
$ws.Range("O3").Value = [double]"0.9191911494312409"
$ws.Range("P3").Value = [double]"0.9191911494312409"
$ws.Range("Q3").Value = [double]"0.6877315928816666"
$ws.Range("R3").Value = [double]"6.189584335935"
$ws.Range("S3").Value = [double]"0.002014630939631782"
$ws.Range("T3").Value = [double]"0.002014630939631783"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"0.040495"
$ws.Range("H4").Value = [double]"0.121485"
$ws.Range("I4").Value = [double]"0.002191743187342868"
$ws.Range("J4").Value = [double]"0.002191743187342869"
$ws.Range("M4").Value = [double]"1.452542333333333"
$ws.Range("N4").Value = [double]"4.357627"
$ws.Range("O4").Value = [double]"0.07861710738141615"
$ws.Range("P4").Value = [double]"0.07861710738141615"
$ws.Range("Q4").Value = [double]"0.05882070178833333"
$ws.Range("R4").Value = [double]"0.529386316095"
$ws.Range("S4").Value = [double]"0.0001723085095118216"
$ws.Range("T4").Value = [double]"0.0001723085095118216"
$ws.Range("I5").Value = [double]"0.9191911494312409"
$ws.Range("J5").Value = [double]"0.9191911494312409"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"0.040495"
$ws.Range("N5").Value = [double]"0.121485"
$ws.Range("O5").Value = [double]"0.002191743187342868"
$ws.Range("P5").Value = [double]"0.002191743187342869"
$ws.Range("Q5").Value = [double]"0.6877315928816666"
$ws.Range("R5").Value = [double]"6.189584335935"
$ws.Range("S5").Value = [double]"0.002014630939631782"
$ws.Range("T5").Value = [double]"0.002014630939631783"
$ws.Range("I6").Value = [double]"0.9191911494312409"
$ws.Range("J6").Value = [double]"0.9191911494312409"
$ws.Range("O6").Value = [double]"0.9191911494312409"
$ws.Range("P6").Value = [double]"0.9191911494312409"
$ws.Range("S6").Value = [double]"0.8449123691927258"
$ws.Range("T6").Value = [double]"0.8449123691927258"
$ws.Range("I7").Value = [double]"0.9191911494312409"
$ws.Range("J7").Value = [double]"0.9191911494312409"
$ws.Range("M7").Value = [double]"1.452542333333333"
$ws.Range("N7").Value = [double]"4.357627"
$ws.Range("O7").Value = [double]"0.07861710738141615"
$ws.Range("P7").Value = [double]"0.07861710738141615"
$ws.Range("S7").Value = [double]"0.07226414929888321"
$ws.Range("T7").Value = [double]"0.07226414929888321"
$ws.Range("G8").Value = [double]"1.452542333333333"
$ws.Range("H8").Value = [double]"4.357627"
$ws.Range("I8").Value = [double]"0.07861710738141615"
$ws.Range("J8").Value = [double]"0.07861710738141615"
$ws.Range("K8").Value = [double]"3"
$ws.Range("L8").Value = [double]"1"
$ws.Range("M8").Value = [double]"0.040495"
$ws.Range("N8").Value = [double]"0.121485"
$ws.Range("O8").Value = [double]"0.002191743187342868"
$ws.Range("P8").Value = [double]"0.002191743187342869"
$ws.Range("Q8").Value = [double]"0.05882070178833333"
$ws.Range("R8").Value = [double]"0.529386316095"
$ws.Range("S8").Value = [double]"0.0001723085095118216"
$ws.Range("T8").Value = [double]"0.0001723085095118216"
$ws.Range("G9").Value = [double]"1.452542333333333"
$ws.Range("H9").Value = [double]"4.357627"
$ws.Range("I9").Value = [double]"0.07861710738141615"
$ws.Range("J9").Value = [double]"0.07861710738141615"
$ws.Range("O9").Value = [double]"0.9191911494312409"
$ws.Range("P9").Value = [double]"0.9191911494312409"
$ws.Range("S9").Value = [double]"0.07226414929888321"
$ws.Range("T9").Value = [double]"0.07226414929888321"
$ws.Range("G10").Value = [double]"1.452542333333333"
$ws.Range("H10").Value = [double]"4.357627"
$ws.Range("I10").Value = [double]"0.07861710738141615"
$ws.Range("J10").Value = [double]"0.07861710738141615"
$ws.Range("M10").Value = [double]"1.452542333333333"
$ws.Range("N10").Value = [double]"4.357627"
$ws.Range("O10").Value = [double]"0.07861710738141615"
$ws.Range("P10").Value = [double]"0.07861710738141615"
$ws.Range("R10").Value = [double]"18.988913071129"
$ws.Range("S10").Value = [double]"0.006180649573021118"
$ws.Range("T10").Value = [double]"0.006180649573021118"

Write-Host "applied 101 cell updates"
